$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had a title row ("Tabela 1") merged across A1:E1 at the
# very top, with the real header row ("Id", "Feature", ...) on row 2 and the
# data starting on row 3. The edit removes that title row entirely so the
# header row becomes row 1 and everything below shifts up by one.
$ws.Rows.Item(1).Delete()
